$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2024-01-05 12:39:36", 13, 9, 0, 1, 5, 3, 4, 0.002, 0.01, 0.001, 100, 500, 10, 7, 0.6923076923076923, "Data/bombay1.xlsx")
    ,@("2024-01-07 23:24:08", 14, 14, 3, 4, 5, 0, 2, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1, "Data/bombay1.xlsx")
    ,@("2024-01-08 01:14:48", 8, 8, 2, 2, 0, 2, 2, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1, "Data/bombay1.xlsx")
    ,@("2024-01-08 11:40:37", 4, 4, 0, 1, 0, 2, 1, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1, "Data/bombay1.xlsx")
    ,@("2024-01-08 12:15:19", 6, 6, 0, 0, 3, 1, 2, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1, "Data/bombay1.xlsx")
    ,@("2024-01-08 19:04:39", 11, 11, 2, 4, 1, 3, 1, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1, "Data/bombay1.xlsx")
    ,@("2024-01-09 00:56:00", 7, 8, 3, 3, 1, 1, 0, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1.142857142857143, "Data/bombay1.xlsx")
    ,@("2024-01-09 01:12:42", 18, 18, 1, 3, 8, 6, 0, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1, "Data/bombay1.xlsx")
    ,@("2024-01-09 01:44:11", 5, 8, 2, 1, 1, 4, 0, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1.6, "Data/bombay1.xlsx")
    ,@("2024-01-09 02:00:12", 1, 1, 0, 1, 0, 0, 0, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1, "Data/lighting.xlsx")
    ,@("2024-01-09 21:16:15", 17, 18, 1, 8, 3, 6, 0, 0.001, 0.05, 0.003, 100, 500, 10, 5, 1.058823529411765, "Data/bombay1.xlsx")
)

$startRow = 63
$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = [string]$row[0]
    for ($ci = 1; $ci -le 15; $ci++) {
        $ws.Cells.Item($r, $ci + 1).Value = $row[$ci]
    }
    $ws.Cells.Item($r, 17).Value = [string]$row[16]
    $r++
}
